# Generate Report for Handback
# Update the timestamp cells that record when the handoff/handback XLIFF
# files were generated, reflecting a re-run of the report generation.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-07 09:37:19"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-07 09:36:59"
$wsZhCn.Range("K2").Value = "2016-09-07 09:38:20"

# de-de sheet: "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-07 09:38:40"
